$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.369.58"
$ws.Range("E2").Value = "  -0.91%  "

$ws.Range("D3").Value = "'3.585.73"
$ws.Range("E3").Value = "  -1.98%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'577.41"
$ws.Range("E5").Value = "  -3.48%  "

$ws.Range("D6").Value = "'191.32"
$ws.Range("E6").Value = "  +0.57%  "

$ws.Range("D7").Value = "'3.577.36"
$ws.Range("E7").Value = "  -2.07%  "

$ws.Range("D8").Value = "'0.617"
$ws.Range("E8").Value = "  -0.54%  "

$ws.Range("E9").Value = "  +0.26%  "

$ws.Range("D10").Value = "'0.677"
$ws.Range("E10").Value = "  -3.39%  "

$ws.Range("D11").Value = "'0.150"
$ws.Range("E11").Value = "  -1.81%  "

$ws.Range("D12").Value = "'54.48"
$ws.Range("E12").Value = "  -5.06%  "

$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("E14").Value = "  -2.61%  "

$ws.Range("D15").Value = "'4.161.52"
$ws.Range("E15").Value = "  -1.77%  "

$ws.Range("D16").Value = "'3.582.58"
$ws.Range("E16").Value = "  -2.07%  "

$ws.Range("E17").Value = "  -0.93%  "

$ws.Range("D18").Value = "'12.31"
$ws.Range("E18").Value = "  -1.05%  "

$ws.Range("D19").Value = "'67.257.76"
$ws.Range("E19").Value = "  -0.64%  "

$ws.Range("E20").Value = "  -2.84%  "

$ws.Range("E21").Value = "  -3.62%  "

$ws.Range("D22").Value = "'402.55"
$ws.Range("E22").Value = "  +0.32%  "

$ws.Range("D23").Value = "'13.33"
$ws.Range("E23").Value = "  +20.21%  "

$ws.Range("E24").Value = "  -4.17%  "

$ws.Range("D25").Value = "'85.63"
$ws.Range("E25").Value = "  -2.31%  "

$ws.Range("E26").Value = "  -0.57%  "

$ws.Range("E27").Value = "  +0.32%  "

$ws.Range("E28").Value = "  +1.07%  "

$ws.Range("D29").Value = "'3.82"
$ws.Range("E29").Value = "  +4.64%  "

$ws.Range("D30").Value = "'8.18"
$ws.Range("E30").Value = "  +12.45%  "

$ws.Range("D31").Value = "'9.12"
$ws.Range("E31").Value = "  -1.79%  "

$ws.Range("D32").Value = "'31.22"
$ws.Range("E32").Value = "  -1.97%  "

$ws.Range("D33").Value = "'670.10"
$ws.Range("E33").Value = "  +10.46%  "

$ws.Range("D34").Value = "'12.20"
$ws.Range("E34").Value = "  -1.01%  "

$ws.Range("E35").Value = "  -0.47%  "

$ws.Range("D36").Value = "'64.01"
$ws.Range("E36").Value = "  -2.90%  "

$ws.Range("D37").Value = "'42.68"
$ws.Range("E37").Value = "  -3.88%  "

$ws.Range("E38").Value = "  +7.30%  "

$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.08%  "

$ws.Range("D40").Value = "'0.0₃0790"
$ws.Range("E40").Value = "  +2.42%  "

$ws.Range("E41").Value = "  +13.59%  "

$ws.Range("D42").Value = "'3.14"
$ws.Range("E42").Value = "  +8.39%  "

$ws.Range("E43").Value = "  -0.61%  "

$ws.Range("D44").Value = "'3.157.02"
$ws.Range("E44").Value = "  +13.70%  "

$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("D46").Value = "'0.0419"
$ws.Range("E46").Value = "  -1.79%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.13"
$ws.Range("E47").Value = "  -1.05%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.131"
$ws.Range("E48").Value = "  -3.20%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'143.46"
$ws.Range("E49").Value = "  +0.25%  "

$ws.Range("D50").Value = "'8.68"
$ws.Range("E50").Value = "  -0.83%  "

$ws.Range("E51").Value = "  -3.02%  "
